# Daily auto-update of the price table: a new row for the latest date is
# inserted at the top (row 2, just below the header), and every existing
# data row shifts down by one. The newest row carries the same price
# values (783.5 / 1112 / 3610) already used throughout the sheet, and the
# date is one day later than the row that used to be on top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 56   # previous last row (55) + 1 new row

# Shift all existing data rows (old row 2..55) down by one (new row 3..56),
# walking bottom-up so we never overwrite a source row before reading it.
for ($r = $lastRow; $r -ge 3; $r--) {
    $src = $r - 1

    $dateText = $ws.Cells.Item($src, 1).Text()
    $ws.Cells.Item($r, 1).Value = "'" + $dateText
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($src, 2).Value()
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($src, 3).Value()
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($src, 4).Value()
}

# New top data row (row 2): latest date, same reference prices.
$ws.Cells.Item(2, 1).Value = "'2026-01-14"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# The apostrophe prefix above forces the date cells to stay plain text
# (matching the rest of the column) instead of being auto-converted to
# date serial numbers; strip the resulting "quote prefix" formatting so
# the cells remain unstyled, just like the rest of the table.
$ws.Range("A2:A" + $lastRow).ClearFormats()
